$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 21 (pushes old rows 21.. down to 23..)
$ws.Rows("21:22").Insert()

# Row 21: new weekly entry - Murcott / Primera
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44469
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100102
$ws.Range("H21").Value = "Cítricos"
$ws.Range("I21").Value = 100102004
$ws.Range("J21").Value = "Mandarina"
$ws.Range("K21").Value = "Murcott"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 5500
$ws.Range("O21").Value = 6000
$ws.Range("P21").Value = 5750
$ws.Range("Q21").Value = "`$/bandeja 10 kilos"
$ws.Range("R21").Value = "Provincia de Limarí"
$ws.Range("S21").Value = 575
$ws.Range("T21").Value = 10

# Row 22: new weekly entry - Murcott / Segunda
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44469
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100102
$ws.Range("H22").Value = "Cítricos"
$ws.Range("I22").Value = 100102004
$ws.Range("J22").Value = "Mandarina"
$ws.Range("K22").Value = "Murcott"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 5000
$ws.Range("O22").Value = 5000
$ws.Range("P22").Value = 5000
$ws.Range("Q22").Value = "`$/bandeja 10 kilos"
$ws.Range("R22").Value = "Provincia de Limarí"
$ws.Range("S22").Value = 500
$ws.Range("T22").Value = 10
